$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.885.13'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.869.55'
$ws.Range('E3').Value = '  +3.02%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '600.77'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '162.75'
$ws.Range('E6').Value = '  -2.74%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.867.57'
$ws.Range('E7').Value = '  +2.99%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.530'
$ws.Range('E9').Value = '  -1.81%  '
$ws.Range('E10').Value = '  -1.38%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.30'
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '36.93'
$ws.Range('E13').Value = '  -2.94%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000243'
$ws.Range('E14').Value = '  -2.19%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.516.01'
$ws.Range('E15').Value = '  +3.00%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.879.79'
$ws.Range('E16').Value = '  +3.43%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '69.057.35'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('E18').Value = '  +2.54%  '
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.10'
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.34'
$ws.Range('E21').Value = '  +1.97%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '484.59'
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('E23').Value = '  -1.40%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.0000162'
$ws.Range('E24').Value = '  +6.01%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '84.02'
$ws.Range('E25').Value = '  -1.11%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.23'
$ws.Range('E26').Value = '  -3.00%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.08'
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.94'
$ws.Range('E29').Value = '  -1.57%  '
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('B31').Value = 'WrappedeETH'
$ws.Range('C31').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.020.69'
$ws.Range('E31').Value = '  +3.03%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.90'
$ws.Range('E32').Value = '  -3.18%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '32.33'
$ws.Range('E33').Value = '  +2.32%  '
$ws.Range('E34').Value = '  -4.18%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.816.68'
$ws.Range('E35').Value = '  +3.48%  '
$ws.Range('E36').Value = '  -1.51%  '
$ws.Range('E37').Value = '  +1.35%  '
$ws.Range('E38').Value = '  +1.66%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.88'
$ws.Range('E39').Value = '  -1.84%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.318'
$ws.Range('E41').Value = '  -2.74%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '436.22'
$ws.Range('E42').Value = '  +1.76%  '
$ws.Range('E43').Value = '  -2.78%  '
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '143.44'
$ws.Range('E48').Value = '  +1.35%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.842.32'
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '26.02'
$ws.Range('E50').Value = '  +13.14%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0357'
$ws.Range('E51').Value = '  +1.02%  '
